$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds text-formatted numbers (e.g. "63.012.27", "1.00", "0.608").
# Excel's COM layer auto-coerces plain numeric-looking strings assigned via
# .Value into real numbers, which would silently change the cell type from
# the original inline/shared string. Force the whole column to Text format
# before writing, then drop back to the Normal style afterwards so no
# leftover style index is left on the cells (matches the source workbook,
# which has no style attribute on these data cells).
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "63.209.86"
$ws.Range("E2").Value = "  -3.39%  "

$ws.Range("D3").Value = "3.296.44"
$ws.Range("E3").Value = "  -5.47%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "543.96"
$ws.Range("E5").Value = "  -1.66%  "

$ws.Range("D6").Value = "171.17"
$ws.Range("E6").Value = "  -4.24%  "

$ws.Range("D7").Value = "0.609"
$ws.Range("E7").Value = "  -4.98%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").Value = "3.288.52"
$ws.Range("E9").Value = "  -5.54%  "

$ws.Range("D10").Value = "0.608"
$ws.Range("E10").Value = "  -3.72%  "

$ws.Range("E11").Value = "  -1.25%  "

$ws.Range("D12").Value = "52.36"
$ws.Range("E12").Value = "  -2.36%  "

$ws.Range("D13").Value = "0.0000263"
$ws.Range("E13").Value = "  -2.73%  "

$ws.Range("D14").Value = "8.83"
$ws.Range("E14").Value = "  -4.19%  "

$ws.Range("D15").Value = "3.813.74"
$ws.Range("E15").Value = "  -5.66%  "

$ws.Range("D16").Value = "17.94"
$ws.Range("E16").Value = "  -2.51%  "

$ws.Range("E17").Value = "  -3.75%  "

$ws.Range("D18").Value = "3.294.43"
$ws.Range("E18").Value = "  -5.47%  "

$ws.Range("D19").Value = "11.59"
$ws.Range("E19").Value = "  -4.51%  "

$ws.Range("D20").Value = "63.024.70"
$ws.Range("E20").Value = "  -3.67%  "

$ws.Range("D21").Value = "0.962"
$ws.Range("E21").Value = "  -3.26%  "

$ws.Range("D22").Value = "422.88"
$ws.Range("E22").Value = "  +2.55%  "

$ws.Range("D23").Value = "4.45"
$ws.Range("E23").Value = "  +8.45%  "

$ws.Range("D24").Value = "4.02"
$ws.Range("E24").Value = "  -0.40%  "

$ws.Range("D25").Value = "13.22"
$ws.Range("E25").Value = "  +3.81%  "

$ws.Range("D26").Value = "82.66"
$ws.Range("E26").Value = "  -3.84%  "

$ws.Range("D27").Value = "10.54"
$ws.Range("E27").Value = "  -2.24%  "

$ws.Range("D28").Value = "2.71"
$ws.Range("E28").Value = "  -4.86%  "

$ws.Range("D29").Value = "8.57"
$ws.Range("E29").Value = "  -5.22%  "

$ws.Range("D30").Value = "28.97"
$ws.Range("E30").Value = "  -4.33%  "

$ws.Range("D31").Value = "6.32"
$ws.Range("E31").Value = "  -2.38%  "

$ws.Range("D32").Value = "11.27"
$ws.Range("E32").Value = "  -3.14%  "

$ws.Range("D33").Value = "571.11"
$ws.Range("E33").Value = "  -6.36%  "

$ws.Range("D34").Value = "0.106"
$ws.Range("E34").Value = "  -3.43%  "

$ws.Range("D35").Value = "57.83"
$ws.Range("E35").Value = "  -2.80%  "

$ws.Range("E36").Value = "  +0.02%  "

$ws.Range("E37").Value = "  -1.83%  "

$ws.Range("E38").Value = "  +4.04%  "

$ws.Range("D39").Value = "34.88"
$ws.Range("E39").Value = "  -5.86%  "

$ws.Range("D40").Value = "0.0₃0735"
$ws.Range("E40").Value = "  -6.58%  "

# Rows 41/42 swap: Maker and TheGraph trade places (with refreshed price
# and volume figures), so write both rows' B/C/D/E together here.
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").Value = "0.362"
$ws.Range("E41").Value = "  -4.62%  "

$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "3.104.56"
$ws.Range("E42").Value = "  -8.12%  "

$ws.Range("D43").Value = "0.998"
$ws.Range("E43").Value = "  -0.18%  "

$ws.Range("D44").Value = "2.75"
$ws.Range("E44").Value = "  -3.10%  "

$ws.Range("E45").Value = "  -0.63%  "

$ws.Range("D46").Value = "0.0399"
$ws.Range("E46").Value = "  -3.53%  "

$ws.Range("D47").Value = "2.41"
$ws.Range("E47").Value = "  -4.44%  "

$ws.Range("D48").Value = "0.128"
$ws.Range("E48").Value = "  -3.65%  "

$ws.Range("E49").Value = "  -5.92%  "

$ws.Range("D50").Value = "131.94"
$ws.Range("E50").Value = "  -3.67%  "

$ws.Range("D51").Value = "8.01"
$ws.Range("E51").Value = "  -4.77%  "

# Restore the default style on the whole column so no stray style index
# (from the temporary Text number format above) remains on any D cell.
$dRange.Style = "Normal"
